# Applies the commit "Added a poll and amended two old polls" to the
# federal.xlsx polling workbook.
#
# 1) Amend the two old Leger national (CAN) rows so that the H/I sample
#    size columns, and some of the G vote-share values, are corrected.
# 2) Append a brand new Abacus poll (pollster id 7, date 44208) with one
#    row per Region/Party combination, following the exact same layout
#    used by the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Amend the two existing Leger "CAN" (national) poll rows.
# ---------------------------------------------------------------------

# Row 38-43: Leger poll #1 (2021-01-?? / serial 44200), region = CAN
$ws.Range("G38").Value = 35
$ws.Range("H38").Value = 1238
$ws.Range("I38").Value = 1242

$ws.Range("G39").Value = 30
$ws.Range("H39").Value = 1238
$ws.Range("I39").Value = 1242

$ws.Range("G40").Value = 20
$ws.Range("H40").Value = 1238
$ws.Range("I40").Value = 1242

$ws.Range("G41").Value = 7
$ws.Range("H41").Value = 1238
$ws.Range("I41").Value = 1242

$ws.Range("H42").Value = 1238
$ws.Range("I42").Value = 1242

$ws.Range("H43").Value = 1238
$ws.Range("I43").Value = 1242

# Row 80-85: Leger poll #2 (serial 44179), region = CAN
$ws.Range("G80").Value = 35
$ws.Range("H80").Value = 1223
$ws.Range("I80").Value = 1248

$ws.Range("G81").Value = 29
$ws.Range("H81").Value = 1223
$ws.Range("I81").Value = 1248

$ws.Range("G82").Value = 23
$ws.Range("H82").Value = 1223
$ws.Range("I82").Value = 1248

$ws.Range("G83").Value = 7
$ws.Range("H83").Value = 1223
$ws.Range("I83").Value = 1248

$ws.Range("G84").Value = 5
$ws.Range("H84").Value = 1223
$ws.Range("I84").Value = 1248

$ws.Range("H85").Value = 1223
$ws.Range("I85").Value = 1248

# ---------------------------------------------------------------------
# 2) Append a new Abacus poll (id 7, date serial 44208) starting at
#    row 190, one row per Region x Party, matching the existing layout
#    (columns A=ID, B=Pollster, C=Date, D=Region, F=Party, G=Vote%).
# ---------------------------------------------------------------------

$pollId = 7
$pollster = "Abacus"
$pollDate = 44208

$newRows = @(
    @('ATL','LIB',38),
    @('ATL','CON',29),
    @('ATL','NDP',17),
    @('ATL','BQ',$null),
    @('ATL','GRN',9),
    @('ATL','PPC',5),
    @('ATL','OTH',1),
    @('QC','LIB',37),
    @('QC','CON',14),
    @('QC','NDP',9),
    @('QC','BQ',35),
    @('QC','GRN',5),
    @('QC','PPC',0),
    @('QC','OTH',1),
    @('ON','LIB',42),
    @('ON','CON',32),
    @('ON','NDP',17),
    @('ON','BQ',$null),
    @('ON','GRN',7),
    @('ON','PPC',1),
    @('ON','OTH',0),
    @('MB/SK','LIB',19),
    @('MB/SK','CON',46),
    @('MB/SK','NDP',24),
    @('MB/SK','BQ',$null),
    @('MB/SK','GRN',3),
    @('MB/SK','PPC',5),
    @('MB/SK','OTH',4),
    @('AB','LIB',21),
    @('AB','CON',54),
    @('AB','NDP',18),
    @('AB','BQ',$null),
    @('AB','GRN',2),
    @('AB','PPC',2),
    @('AB','OTH',3),
    @('BC','LIB',29),
    @('BC','CON',30),
    @('BC','NDP',29),
    @('BC','BQ',$null),
    @('BC','GRN',11),
    @('BC','PPC',0),
    @('BC','OTH',0),
    @('CAN','LIB',35),
    @('CAN','CON',31),
    @('CAN','NDP',17),
    @('CAN','BQ',8),
    @('CAN','GRN',6),
    @('CAN','PPC',1),
    @('CAN','OTH',1)
)

$row = 190
foreach ($entry in $newRows) {
    $region = $entry[0]
    $party = $entry[1]
    $vote = $entry[2]

    $ws.Cells.Item($row, 1).Value = $pollId
    $ws.Cells.Item($row, 2).Value = $pollster
    $ws.Cells.Item($row, 3).Value = $pollDate
    $ws.Cells.Item($row, 3).NumberFormat = "yyyy/mm/dd"
    $ws.Cells.Item($row, 4).Value = $region
    $ws.Cells.Item($row, 6).Value = $party
    if ($vote -ne $null) {
        $ws.Cells.Item($row, 7).Value = $vote
    }

    $row = $row + 1
}

# ---------------------------------------------------------------------
# Cosmetic: reflect the final selection used in the original edit
# (the author ended up with H85:I85 selected on the frozen pane).
# ---------------------------------------------------------------------
$ws.Range("H85:I85").Select() | Out-Null
